# Commit: "Updated notebook, reran simulation"
#
# Two new simulated textures ("Holden" and "Rizzie Spiral") were added to
# the table right after the "Spiral5" row, pushing every following row
# down by two rows and renumbering the index column (A) accordingly.
# "Thomas Hex" was also renamed to "Matthies Hex" (same data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1. Shift the existing rows 4..29 down to rows 6..31 (columns B..T only -
#    column A's running index is rebuilt afterwards). Go bottom-up so
#    earlier rows are not overwritten before they are copied.
# -----------------------------------------------------------------------
for ($r = 29; $r -ge 4; $r--) {
    $src = $ws.Range("B" + $r + ":T" + $r)
    $dst = $ws.Range("B" + ($r + 2) + ":T" + ($r + 2))
    $dst.Value = $src.Value()
}

# -----------------------------------------------------------------------
# 2. Write the freshly-simulated data for the two new rows.
# -----------------------------------------------------------------------
$ws.Cells.Item(4, 2).Value = "Holden"
$row4vals = @(0.9986984847922966,1.001577596349724,1.000591601450591,1.000591601450591,0.9976336060975058,1.000591601450591,0.9976336060975058,1.000591601450591,1.000591601450591,1.000591601450591,0.9991126037740485,0.9991126037740485,0.9989745641134645,0.9996056029995627,0.9996056029995627,0.9998521026123199,0.9998521026123199,0.9999474152652167)
$col = 3
foreach ($v in $row4vals) {
    $ws.Cells.Item(4, $col).Value = $v
    $col = $col + 1
}

$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"
$row5vals = @(1.015252498150222,0.9815121182034636,0.99306704431791,0.99306704431791,1.027731828751026,0.99306704431791,1.027731828751026,0.99306704431791,0.99306704431791,0.99306704431791,1.010399436534468,1.010399436534468,1.012017123739719,1.004621972462282,1.004621972462282,1.001733240426189,1.001733240426189,1.00061626300974)
$col = 3
foreach ($v in $row5vals) {
    $ws.Cells.Item(5, $col).Value = $v
    $col = $col + 1
}

# -----------------------------------------------------------------------
# 3. Rebuild the sequential index column (A) for the whole table body.
# -----------------------------------------------------------------------
for ($r = 3; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# -----------------------------------------------------------------------
# 4. The table grew by two rows, so the two new bottom rows (30:31) need
#    the bold/bordered "index" formatting that column A uses throughout;
#    every other row already carries it.
# -----------------------------------------------------------------------
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# -----------------------------------------------------------------------
# 5. Rename "Thomas Hex" -> "Matthies Hex" (now at row 11, since it used
#    to be row 9 before the two new rows were inserted above it).
# -----------------------------------------------------------------------
$ws.Cells.Item(11, 2).Value = "Matthies Hex"
